$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column in front of the existing "Unlocked" column (column E),
#    pushing Mountains/Icy planes/Swamps/Prairie-Savannah/Forest?/Comments one
#    column to the right (E:K -> F:L). Excel copies the left neighbour's
#    (column D) number format / borders onto the freshly inserted column, which
#    is what the target file also does for the data rows.
$ws.Columns.Item(5).Insert()

# 2. The new column is a second "Unlocked" style column that records whether a
#    resource unlocks the "Recovery" category - give it a header.
$hdr = $ws.Range("E1")
$hdr.Value = "Recovery"

# 3. Style the new header cell like the other header cells (bold font already
#    carried over from the insert) but flag it with its own accent fill so it
#    stands out from the rest of the header row.
$hdr.Interior.ThemeColor = 5
$hdr.Interior.TintAndShade = 0.8

# 4. Give the new column its own width (content was short: "Recovery").
$ws.Range("E1").ColumnWidth = 10.67

# 5. The autofilter / filter database range needs to grow from K37 to L37 to
#    include the newly inserted column. Re-applying AutoFilter on the wider
#    range (after turning the stale one off) refreshes the stored ref.
$ws.AutoFilterMode = $false
$ws.Range("A1:L37").AutoFilter()

$filterName = $wb.Names.Item("Tabelle1!_FilterDatabase")
$filterName.RefersTo = "=Tabelle1!`$A`$1:`$L`$37"

# 6. Leave the cursor roughly where editing finished.
$ws.Range("F4").Select()
